$wb = $excel.ActiveWorkbook

# Rename the "摘要" sheet to "交易別"
$ws4 = $wb.Worksheets.Item("摘要")
$ws4.Name = "交易別"

# On the "JsonField" sheet, insert a new row before row 35 describing the
# new "Memo" / "摘要" field (NVARCHAR2, length 60).
$ws3 = $wb.Worksheets.Item("JsonField")
$ws3.Rows.Item(35).Insert()

$ws3.Cells.Item(35, 1).Value = 33
$ws3.Cells.Item(35, 2).Value = "Memo"
$ws3.Cells.Item(35, 3).Value = "摘要"
$ws3.Cells.Item(35, 4).Value = "NVARCHAR2"
$ws3.Cells.Item(35, 5).Value = 60

# The row-insert copies row 34's formatting down, which leaves D35 on the
# wrong style and creates a stray formatted H35. Line both back up with
# the rest of the "交易別" field block (row 36 onward).
$ws3.Cells.Item(36, 4).Copy()
$ws3.Cells.Item(35, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws3.Cells.Item(35, 8).Clear()

$ws3.Activate()
